# Weekly update for "Fruta / hortaliza" (Mango, Vega Modelo de Temuco):
# a new week's record is inserted as row 208, pushing the existing
# rows 208-238 down to 209-239 (their own data is preserved as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 208; Excel shifts rows 208:238 down to 209:239.
$ws.Rows("208:208").Insert()

# Populate the newly inserted row with this week's reading.
$ws.Cells.Item(208, 1).Value  = 10
$ws.Cells.Item(208, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(208, 3).Value  = "La Araucanía"
$ws.Cells.Item(208, 4).Value  = 44504
$ws.Cells.Item(208, 5).Value  = 9
$ws.Cells.Item(208, 6).Value  = "Fruta"
$ws.Cells.Item(208, 7).Value  = 100108
$ws.Cells.Item(208, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(208, 9).Value  = 100108002
$ws.Cells.Item(208, 10).Value = "Mango"
$ws.Cells.Item(208, 11).Value = "Sin especificar"
$ws.Cells.Item(208, 12).Value = "Primera"
$ws.Cells.Item(208, 13).Value = 2000
$ws.Cells.Item(208, 14).Value = 8000
$ws.Cells.Item(208, 15).Value = 8000
$ws.Cells.Item(208, 16).Value = 8000
$ws.Cells.Item(208, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(208, 18).Value = "Perú"
$ws.Cells.Item(208, 19).Value = 2000
$ws.Cells.Item(208, 20).Value = 4
